$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '25.840.74'
$ws.Range('E2').Value = '  -1.22%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.633.04'
$ws.Range('E3').Value = '  -1.35%  '
$ws.Range('E4').Value = '  -0.39%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '214.83'
$ws.Range('E5').Value = '  -0.57%  '
$ws.Range('E6').Value = '  -1.79%  '
$ws.Range('E7').Value = '  -0.39%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2565'
$ws.Range('E8').Value = '  -0.97%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06400'
$ws.Range('E9').Value = '  -0.26%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.57'
$ws.Range('E10').Value = '  -2.01%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07684'
$ws.Range('E11').Value = '  -1.91%  '
$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.640.62'
$ws.Range('E12').Value = '  -0.94%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.237'
$ws.Range('E13').Value = '  -1.27%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.859.12'
$ws.Range('E14').Value = '  -1.29%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.5422'
$ws.Range('E15').Value = '  -2.18%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0₅7919'
$ws.Range('E16').Value = '  -1.34%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '63.37'
$ws.Range('E17').Value = '  -1.26%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '25.850.97'
$ws.Range('E18').Value = '  -1.32%  '
$ws.Range('E19').Value = '  -0.35%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '201.54'
$ws.Range('E20').Value = '  -4.09%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.320'
$ws.Range('E21').Value = '  -2.31%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.920'
$ws.Range('E22').Value = '  -1.50%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.967'
$ws.Range('E23').Value = '  -1.26%  '
$ws.Range('E24').Value = '  -0.30%  '
$ws.Range('E25').Value = '  +10.47%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '140.83'
$ws.Range('E26').Value = '  -1.85%  '
$ws.Range('E27').Value = '  -2.67%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '15.66'
$ws.Range('E28').Value = '  -0.99%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '6.690'
$ws.Range('E29').Value = '  -4.29%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.238'
$ws.Range('E30').Value = '  -0.73%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.04971'
$ws.Range('E31').Value = '  -2.85%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.254'
$ws.Range('E32').Value = '  -3.09%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.172'
$ws.Range('E33').Value = '  -1.73%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.535'
$ws.Range('E34').Value = '  -1.89%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.363'
$ws.Range('E35').Value = '  -0.44%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.165.35'
$ws.Range('E36').Value = '  +1.21%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.8911'
$ws.Range('E37').Value = '  -4.11%  '
$ws.Range('E38').Value = '  -4.72%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.5592'
$ws.Range('E39').Value = '  -2.05%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.01557'
$ws.Range('E40').Value = '  -2.33%  '
$ws.Range('E41').Value = '  -0.35%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.544'
$ws.Range('E42').Value = '  -0.62%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.673'
$ws.Range('E43').Value = '  +0.44%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.8060'
$ws.Range('E44').Value = '  -3.64%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '99.23'
$ws.Range('E45').Value = '  -1.27%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.771.47'
$ws.Range('E46').Value = '  -1.25%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0₈114'
$ws.Range('E47').Value = '  -0.65%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.4516'
$ws.Range('E48').Value = '  -0.69%  '
$ws.Range('E49').Value = '  -0.11%  '
$ws.Range('E50').Value = '  -2.07%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.05074'
$ws.Range('E51').Value = '  +0.56%  '
